# "Generate Report for Handback" -- refresh the handback status/report values
# across the Overview, zh-cn and de-de sheets, and widen the Status / Error
# Detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Widen columns E and F to fit the longer status text (best-effort match of
# the target stored width, which is snapped to this engine's column-width
# grid).
$wsOverview.Columns.Item(5).ColumnWidth = 29.14437166849777
$wsOverview.Columns.Item(6).ColumnWidth = 29.14437166849777

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-10-20 08:58:32"
$wsZhCn.Range("P2").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719813028965

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-10-20 08:58:50"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14437166849777
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719813028965
